# Act 2 Lilith / Scene 72
#
# Prim's "Sorry you have to do this…" line was split across three runs:
#   "Prim ("  +  "shy worried_slightly"  +  "): Sorry you have to do this…"
# Re-save it as a single contiguous run with the same formatting by doing
# a literal find/replace over the full phrase - Word's Find/Replace
# consolidates the matched text back into one run.
$d = $word.ActiveDocument

$phrase = "Prim (shy worried_slightly): Sorry you have to do this…"

$d.Content.Find.Execute($phrase, $true, $false, $false, $false, $false, $true, 1, $false, $phrase, 2)
